# Update "想去人数" (interested-people count) figures in column F
# for the sheets that hold the full data table: "展览" and "全部类型".
# Both sheets mirror the same rows, so the same updates are applied twice.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1084
    7  = 2438
    11 = 1233
    15 = 1097
    16 = 301
    17 = 312
    19 = 21
    21 = 64
    23 = 151
    24 = 11
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
